$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '37.827.39'
$ws.Range('E2').Value = '  +0.96%  '

$ws.Range('D3').Value = '2.088.59'
$ws.Range('E3').Value = '  +0.82%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '235.11'
$ws.Range('E5').Value = '  +0.05%  '

$ws.Range('E6').Value = '  +0.09%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '58.77'
$ws.Range('E7').Value = '  +2.59%  '

$ws.Range('E9').Value = '  +0.27%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0790'
$ws.Range('E10').Value = '  +1.51%  '

$ws.Range('E11').Value = '  +2.70%  '

$ws.Range('D12').Value = '2.395.15'
$ws.Range('E12').Value = '  +0.90%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '14.77'
$ws.Range('E13').Value = '  +2.82%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '21.29'
$ws.Range('E14').Value = '  +2.92%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.773'
$ws.Range('E15').Value = '  -0.84%  '

$ws.Range('E16').Value = '  +2.35%  '

$ws.Range('D17').Value = '2.085.31'
$ws.Range('E17').Value = '  +1.40%  '

$ws.Range('D18').Value = '37.763.06'
$ws.Range('E18').Value = '  +1.10%  '

$ws.Range('E19').Value = '  +0.05%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '71.31'
$ws.Range('E20').Value = '  +2.44%  '

$ws.Range('D21').Value = '0.0₃0836'
$ws.Range('E21').Value = '  +2.00%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '228.43'
$ws.Range('E22').Value = '  +0.86%  '

$ws.Range('E23').Value = '  -0.10%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.46'
$ws.Range('E24').Value = '  +2.43%  '

$ws.Range('E25').Value = '  -1.20%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '170.39'
$ws.Range('E26').Value = '  +1.14%  '

$ws.Range('E27').Value = '  +4.47%  '

$ws.Range('E28').Value = '  +2.12%  '

$ws.Range('B29').Value = 'EthereumClassic'
$ws.Range('C29').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '19.53'
$ws.Range('E29').Value = '  +2.25%  '

$ws.Range('B30').Value = 'ImmutableX'
$ws.Range('C30').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.40'
$ws.Range('E30').Value = '  -0.60%  '

$ws.Range('E31').Value = '  +2.37%  '

$ws.Range('E32').Value = '  +2.84%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0634'
$ws.Range('E33').Value = '  +2.60%  '

$ws.Range('E34').Value = '  +3.28%  '

$ws.Range('E35').Value = '  +0.99%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '3.46'
$ws.Range('E36').Value = '  +2.76%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.84'
$ws.Range('E37').Value = '  +2.93%  '

$ws.Range('E38').Value = '  -0.02%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.39'
$ws.Range('E39').Value = '  -3.98%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0995'
$ws.Range('E40').Value = '  +4.23%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '98.76'
$ws.Range('E41').Value = '  +1.80%  '

$ws.Range('E42').Value = '  -0.33%  '

$ws.Range('B43').Value = 'VeChain'
$ws.Range('C43').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0216'
$ws.Range('E43').Value = '  +1.40%  '

$ws.Range('B44').Value = 'Maker'
$ws.Range('C44').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D44').Value = '1.464.10'
$ws.Range('E44').Value = '  -1.82%  '

$ws.Range('B45').Value = 'FTXToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '4.34'
$ws.Range('E45').Value = '  +3.78%  '

$ws.Range('E46').Value = '  +0.65%  '

$ws.Range('E47').Value = '  +4.13%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '16.02'
$ws.Range('E48').Value = '  +5.11%  '

$ws.Range('E49').Value = '  +2.59%  '

$ws.Range('E50').Value = '  +2.52%  '

$ws.Range('D51').Value = '2.279.18'
$ws.Range('E51').Value = '  +0.81%  '
